# Add instructions to running tests locally to test_report.
# The "Comments" cell (B18) on the active Bug Report sheet explained the
# Seriousness rating as "critical"; update it to say "major" (matching the
# "Major" seriousness already recorded in B13), and leave the selection on
# that cell as the user would after editing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "Function is used to show price information to the user and that is why Seriousness is determined to be major. Fix should be straightforward to implement."

$ws.Range("B18").Select()
